$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New status entries (Nov 3rd) ---

# Row 77: single wrapped/centered note spanning column B
$ws.Range("B77").Value = "Discussed in the group , got the doubts cleared on concepts finished till now"
$ws.Range("B77").WrapText = $true
$ws.Range("B77").HorizontalAlignment = -4108
$ws.Range("B77").VerticalAlignment = -4108
$ws.Rows.Item(77).RowHeight = 31.5

# Row 78: DONE / PROGRESS / TO-DO style row across B:D
$ws.Range("B78").Value = "Started exploring the codec 2.0 and collected the useful links"
$ws.Range("B78").WrapText = $true
$ws.Range("B78").HorizontalAlignment = -4108
$ws.Range("B78").VerticalAlignment = -4108

$ws.Range("C78").Value = "Studying the collected pdfs on OPENMAX IL "
$ws.Range("D78").Value = "Revision of C-DS-OS concepts"

# Row 79: trailing wrapped/centered note spanning column B
$ws.Range("B79").Value = "Discussed doubts on ffmpeg with the new teammate"
$ws.Range("B79").WrapText = $true
$ws.Range("B79").HorizontalAlignment = -4108
$ws.Range("B79").VerticalAlignment = -4108

# --- View state: zoom in and scroll down to the newly added rows ---
$excel.ActiveWindow.Zoom = 85
[void]$ws.Range("D78").Select()
